$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while writing values so numeric-looking
# strings (e.g. "583.82") are not auto-converted to numbers by COM type
# inference; the source workbook stores every D-column cell as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.134.01'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '2.486.64'
$ws.Range("E3").Value = '  +0.54%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '583.82'
$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("D6").Value = '171.50'
$ws.Range("E6").Value = '  +4.24%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").Value = '2.486.50'
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("E10").Value = '  +1.51%  '

$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("E13").Value = '  -1.73%  '

$ws.Range("D14").Value = '2.960.40'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").Value = '66.926.79'
$ws.Range("E16").Value = '  +0.23%  '

$ws.Range("E17").Value = '  -1.15%  '

$ws.Range("D18").Value = '2.483.12'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = '11.00'
$ws.Range("E19").Value = '  -4.70%  '

$ws.Range("E20").Value = '  -4.97%  '

$ws.Range("D21").Value = '348.62'
$ws.Range("E21").Value = '  -2.60%  '

$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("B24").Value = 'NEARProtocol'
$ws.Range("C24").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D24").Value = '4.23'
$ws.Range("E24").Value = '  -3.59%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '68.37'
$ws.Range("E25").Value = '  -2.79%  '

$ws.Range("D26").Value = '1.80'
$ws.Range("E26").Value = '  -2.08%  '

$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").Value = '0.0₃0905'
$ws.Range("E30").Value = '  -2.42%  '

$ws.Range("D31").Value = '510.93'
$ws.Range("E31").Value = '  +2.60%  '

$ws.Range("E32").Value = '  -3.67%  '

$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -2.27%  '

$ws.Range("E34").Value = '  -3.54%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("E37").Value = '  -6.34%  '

$ws.Range("D38").Value = '18.70'
$ws.Range("E38").Value = '  +0.74%  '

$ws.Range("E39").Value = '  -3.40%  '

$ws.Range("E40").Value = '  -4.83%  '

$ws.Range("E41").Value = '  -1.71%  '

$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("E43").Value = '  -1.49%  '

$ws.Range("E44").Value = '  -2.02%  '

$ws.Range("E45").Value = '  -3.54%  '

$ws.Range("D46").Value = '38.81'
$ws.Range("E46").Value = '  -1.11%  '

$ws.Range("D47").Value = '142.89'
$ws.Range("E47").Value = '  +1.31%  '

$ws.Range("D48").Value = '3.46'
$ws.Range("E48").Value = '  -4.00%  '

$ws.Range("D49").Value = '0.514'
$ws.Range("E49").Value = '  -3.76%  '

$ws.Range("D50").Value = '0.0₆0252'
$ws.Range("E50").Value = '  -3.62%  '

$ws.Range("E51").Value = '  -0.92%  '

# Restore the original (default/General) number format now that the
# text values are committed, so no stray style is left on the cells.
$ws.Range("D2:D51").ClearFormats()
